# Auto-generated Excel COM-interop script
# Applies updated market/profit data values per sheet, matching the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1714.091
$ws.Range("J17").Value = 1971.1364
$ws.Range("L17").Value = 5913.4092
$ws.Range("N17").Value = -6249.4092
$ws.Range("H28").Value = 542
$ws.Range("I28").Value = 452.375
$ws.Range("K28").Value = 452.375
$ws.Range("M28").Value = 32.625
$ws.Range("H62").Value = 22654.6
$ws.Range("J62").Value = 33720.332
$ws.Range("L62").Value = 33720.332
$ws.Range("N62").Value = -34968.332
$ws.Range("H65").Value = 22654.6
$ws.Range("J65").Value = 33720.332
$ws.Range("L65").Value = 168601.66
$ws.Range("N65").Value = -174841.66
$ws.Range("H98").Value = 1863.5646
$ws.Range("I98").Value = 1960.4546
$ws.Range("J98").Value = 1102.2858
$ws.Range("K98").Value = 1960.4546
$ws.Range("L98").Value = 1102.2858
$ws.Range("M98").Value = -462.4546
$ws.Range("N98").Value = -4098.2858
$ws.Range("H99").Value = 1771.3334
$ws.Range("I99").Value = 1747
$ws.Range("J99").Value = 1783.5
$ws.Range("K99").Value = 5241
$ws.Range("L99").Value = 5350.5
$ws.Range("M99").Value = -3743
$ws.Range("N99").Value = -8346.5
$ws.Range("H116").Value = 33125
$ws.Range("I116").Value = 37500
$ws.Range("K116").Value = 37500
$ws.Range("M116").Value = -34058
$ws.Range("H122").Value = 1863.5646
$ws.Range("I122").Value = 1960.4546
$ws.Range("J122").Value = 1102.2858
$ws.Range("K122").Value = 5881.3638
$ws.Range("L122").Value = 3306.8574
$ws.Range("M122").Value = -3431.3638
$ws.Range("N122").Value = -8206.857400000001
$ws.Range("H129").Value = 3109.0715
$ws.Range("I129").Value = 2900.9092
$ws.Range("J129").Value = 3872.3333
$ws.Range("K129").Value = 8702.7276
$ws.Range("L129").Value = 11616.9999
$ws.Range("M129").Value = -3702.7276
$ws.Range("N129").Value = -21616.9999
$ws.Range("H132").Value = 475132.34
$ws.Range("I132").Value = 515143.53
$ws.Range("J132").Value = 8335
$ws.Range("K132").Value = 1545430.59
$ws.Range("L132").Value = 25005
$ws.Range("M132").Value = -1542900.59
$ws.Range("N132").Value = -30065
$ws.Range("H138").Value = 3379.926
$ws.Range("I138").Value = 1688.091
$ws.Range("J138").Value = 4543.0625
$ws.Range("K138").Value = 5064.272999999999
$ws.Range("L138").Value = 13629.1875
$ws.Range("M138").Value = 75.72700000000077
$ws.Range("N138").Value = -23909.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2241.6453
$ws.Range("I2").Value = 1586.8948
$ws.Range("J2").Value = 3278.3333
$ws.Range("K2").Value = 1586.8948
$ws.Range("L2").Value = 3278.3333
$ws.Range("M2").Value = -1473.8948
$ws.Range("N2").Value = -3504.3333
$ws.Range("H32").Value = 20014526
$ws.Range("I32").Value = 26327110
$ws.Range("J32").Value = 24678
$ws.Range("K32").Value = 26327110
$ws.Range("L32").Value = 24678
$ws.Range("M32").Value = -26326823
$ws.Range("N32").Value = -25252
$ws.Range("H45").Value = 6666
$ws.Range("I45").Value = 7999
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 7999
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -7622
$ws.Range("N45").Value = -4754
$ws.Range("H61").Value = 6898.4062
$ws.Range("I61").Value = 3087.6191
$ws.Range("K61").Value = 3087.6191
$ws.Range("M61").Value = -2875.6191
$ws.Range("H74").Value = 3381121.2
$ws.Range("I74").Value = 4311919.5
$ws.Range("K74").Value = 4311919.5
$ws.Range("M74").Value = -4311045.5
$ws.Range("H77").Value = 3381121.2
$ws.Range("I77").Value = 4311919.5
$ws.Range("K77").Value = 21559597.5
$ws.Range("M77").Value = -21555229.5
$ws.Range("H116").Value = 2241.6453
$ws.Range("I116").Value = 1586.8948
$ws.Range("J116").Value = 3278.3333
$ws.Range("K116").Value = 1586.8948
$ws.Range("L116").Value = 3278.3333
$ws.Range("M116").Value = 707.1052
$ws.Range("N116").Value = -7866.3333
$ws.Range("H122").Value = 2529.9375
$ws.Range("I122").Value = 2652.2307
$ws.Range("K122").Value = 7956.6921
$ws.Range("M122").Value = -5506.6921
$ws.Range("H136").Value = 6898.4062
$ws.Range("I136").Value = 3087.6191
$ws.Range("K136").Value = 9262.8573
$ws.Range("M136").Value = -6712.8573

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2241.6453
$ws.Range("I3").Value = 1586.8948
$ws.Range("J3").Value = 3278.3333
$ws.Range("K3").Value = 1586.8948
$ws.Range("L3").Value = 3278.3333
$ws.Range("M3").Value = -1472.8948
$ws.Range("N3").Value = -3506.3333
$ws.Range("H105").Value = 3611.875
$ws.Range("I105").Value = 3317.5925
$ws.Range("J105").Value = 5201
$ws.Range("K105").Value = 3317.5925
$ws.Range("L105").Value = 5201
$ws.Range("M105").Value = -1570.5925
$ws.Range("N105").Value = -8695
$ws.Range("H134").Value = 631864.3
$ws.Range("I134").Value = 885975.5600000001
$ws.Range("J134").Value = 8136.636
$ws.Range("K134").Value = 2657926.68
$ws.Range("L134").Value = 24409.908
$ws.Range("M134").Value = -2655391.68
$ws.Range("N134").Value = -29479.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11478.343
$ws.Range("I31").Value = 1833.2222
$ws.Range("K31").Value = 1833.2222
$ws.Range("M31").Value = -1538.2222
$ws.Range("H34").Value = 11478.343
$ws.Range("I34").Value = 1833.2222
$ws.Range("K34").Value = 1833.2222
$ws.Range("M34").Value = -1631.2222
$ws.Range("H94").Value = 17220
$ws.Range("J94").Value = 2755
$ws.Range("L94").Value = 2755
$ws.Range("N94").Value = -3657
$ws.Range("H122").Value = 1933.1923
$ws.Range("I122").Value = 1684.6818
$ws.Range("K122").Value = 5054.0454
$ws.Range("M122").Value = -2604.0454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2206.923
$ws.Range("I5").Value = 478.8
$ws.Range("J5").Value = 7967.3335
$ws.Range("K5").Value = 1436.4
$ws.Range("L5").Value = 23902.0005
$ws.Range("M5").Value = -1324.4
$ws.Range("N5").Value = -24126.0005
$ws.Range("H12").Value = 62.18182
$ws.Range("I12").Value = 92.59999999999999
$ws.Range("K12").Value = 277.8
$ws.Range("M12").Value = -104.8
$ws.Range("H60").Value = 1201.6111
$ws.Range("I60").Value = 1134.9286
$ws.Range("K60").Value = 3404.7858
$ws.Range("M60").Value = -3153.7858
$ws.Range("H81").Value = 5705.5
$ws.Range("J81").Value = 6631.875
$ws.Range("L81").Value = 19895.625
$ws.Range("N81").Value = -22141.625
$ws.Range("H84").Value = 5705.5
$ws.Range("J84").Value = 6631.875
$ws.Range("L84").Value = 59686.875
$ws.Range("N84").Value = -70918.875
$ws.Range("H86").Value = 1041.7
$ws.Range("J86").Value = 1036.1666
$ws.Range("L86").Value = 3108.4998
$ws.Range("N86").Value = -5480.4998
$ws.Range("H89").Value = 1041.7
$ws.Range("J89").Value = 1036.1666
$ws.Range("L89").Value = 9325.499400000001
$ws.Range("N89").Value = -21181.4994
$ws.Range("H107").Value = 933.5714
$ws.Range("I107").Value = 707.2
$ws.Range("J107").Value = 1499.5
$ws.Range("K107").Value = 2121.6
$ws.Range("L107").Value = 4498.5
$ws.Range("M107").Value = -201.6000000000004
$ws.Range("N107").Value = -8338.5
$ws.Range("H113").Value = 1418.25
$ws.Range("I113").Value = 1257
$ws.Range("J113").Value = 1498.875
$ws.Range("K113").Value = 3771
$ws.Range("L113").Value = 4496.625
$ws.Range("M113").Value = -1601
$ws.Range("N113").Value = -8836.625
$ws.Range("H135").Value = 2206.923
$ws.Range("I135").Value = 478.8
$ws.Range("J135").Value = 7967.3335
$ws.Range("K135").Value = 4309.2
$ws.Range("L135").Value = 71706.0015
$ws.Range("M135").Value = -1774.2
$ws.Range("N135").Value = -76776.0015
$ws.Range("H139").Value = 2015
$ws.Range("J139").Value = 5000
$ws.Range("L139").Value = 15000
$ws.Range("N139").Value = -25280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31262
$ws.Range("H113").Value = 1673
$ws.Range("I113").Value = 1700
$ws.Range("J113").Value = 1619
$ws.Range("K113").Value = 1700
$ws.Range("L113").Value = 1619
$ws.Range("M113").Value = 470
$ws.Range("N113").Value = -5959
$ws.Range("H122").Value = 3273.9688
$ws.Range("I122").Value = 1853.6296
$ws.Range("K122").Value = 5560.8888
$ws.Range("M122").Value = -3110.8888

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1199
$ws.Range("I40").Value = 1199
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1199
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1063
$ws.Range("N40").ClearContents()
$ws.Range("H132").Value = 1085500.1
$ws.Range("I132").Value = 1445002.4
$ws.Range("J132").Value = 6993.625
$ws.Range("K132").Value = 4335007.199999999
$ws.Range("L132").Value = 20980.875
$ws.Range("M132").Value = -4332477.199999999
$ws.Range("N132").Value = -26040.875
$ws.Range("H136").Value = 6285.5625
$ws.Range("I136").Value = 4055.8333
$ws.Range("K136").Value = 12167.4999
$ws.Range("M136").Value = -9617.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H122").Value = 2665.261
$ws.Range("I122").Value = 2121.158
$ws.Range("K122").Value = 6363.474
$ws.Range("M122").Value = -3913.474
$ws.Range("H126").Value = 4334.8125
$ws.Range("I126").Value = 4305.727
$ws.Range("K126").Value = 12917.181
$ws.Range("M126").Value = -10447.181
$ws.Range("H132").Value = 8249844.5
$ws.Range("I132").Value = 653422.9
$ws.Range("K132").Value = 1960268.7
$ws.Range("M132").Value = -1957738.7
$ws.Range("H133").Value = 57115.168
$ws.Range("J133").Value = 57115.168
$ws.Range("L133").Value = 57115.168
$ws.Range("N133").Value = -67235.16800000001
